$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "Trends Status" sheet - update values (rows 2-8)
# ---------------------------------------------------------------------------
$wsTrends = $wb.Worksheets.Item("Trends Status")

$wsTrends.Cells.Item(2,2).Value = 0      # B2 Rapid Decline - Long-term species (no.)
$wsTrends.Cells.Item(2,3).Value = 1      # C2
$wsTrends.Cells.Item(2,4).Value = 0      # D2
$wsTrends.Cells.Item(2,5).Value = 50     # E2

$wsTrends.Cells.Item(3,2).Value = 1      # B3 Decline
$wsTrends.Cells.Item(3,3).Value = 1      # C3 (unchanged)
$wsTrends.Cells.Item(3,4).Value = 33.3   # D3
$wsTrends.Cells.Item(3,5).Value = 50     # E3

$wsTrends.Cells.Item(4,2).Value = 0      # B4 Stable
$wsTrends.Cells.Item(4,3).Value = 0      # C4
$wsTrends.Cells.Item(4,4).Value = 0      # D4
$wsTrends.Cells.Item(4,5).Value = 0      # E4

$wsTrends.Cells.Item(5,2).Value = 0      # B5 Increase
$wsTrends.Cells.Item(5,3).Value = 0      # C5 (unchanged)
$wsTrends.Cells.Item(5,4).Value = 0      # D5
$wsTrends.Cells.Item(5,5).Value = 0      # E5 (unchanged)

$wsTrends.Cells.Item(6,2).Value = 2      # B6 Rapid Increase
$wsTrends.Cells.Item(6,3).Value = 0      # C6
$wsTrends.Cells.Item(6,4).Value = 66.7   # D6
$wsTrends.Cells.Item(6,5).Value = 0      # E6

$wsTrends.Cells.Item(7,2).Value = 22     # B7 Trend Inconclusive
$wsTrends.Cells.Item(7,3).Value = 28     # C7

$wsTrends.Cells.Item(8,2).Value = 313    # B8 Insufficient Data
$wsTrends.Cells.Item(8,3).Value = 308    # C8

# ---------------------------------------------------------------------------
# 2. "Priority Status" sheet - update species counts
# ---------------------------------------------------------------------------
$wsPriority = $wb.Worksheets.Item("Priority Status")

$wsPriority.Cells.Item(2,2).Value = 103  # B2 High
$wsPriority.Cells.Item(3,2).Value = 286  # B3 Moderate
$wsPriority.Cells.Item(4,2).Value = 554  # B4 Low

# ---------------------------------------------------------------------------
# 3. "Species qualification" sheet
# ---------------------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")

$wsQual.Cells.Item(2,1).Value = "SoIB Assessment"  # A2 (was "SoIB 2023 Assessment")
$wsQual.Cells.Item(2,2).Value = 338                # B2

$wsQual.Cells.Item(3,2).Value = 25                 # B3 Long-term Analysis
$wsQual.Cells.Item(3,3).Value = 3                  # C3

$wsQual.Cells.Item(4,2).Value = 30                 # B4 Current Analysis
$wsQual.Cells.Item(4,3).Value = 2                  # C4

# ---------------------------------------------------------------------------
# 4. Rename "High Priority break-up" -> "Interannual update - High Pri"
#    and replace its values with the new interannual-update figures.
# ---------------------------------------------------------------------------
$wsInterannual = $wb.Worksheets.Item("High Priority break-up")
$wsInterannual.Name = "Interannual update - High Pri"

$wsInterannual.Cells.Item(2,2).Value = 72          # B2 Trend New - High Species (no.)
$wsInterannual.Cells.Item(2,3).Value = 69.90000000000001   # C2 High Species (perc.)
$wsInterannual.Cells.Item(2,4).Value = 72          # D2 New High Species (no.)
$wsInterannual.Cells.Item(2,5).Value = 75.8        # E2 New High Species (perc.)

$wsInterannual.Cells.Item(3,2).Value = 31          # B3 IUCN
$wsInterannual.Cells.Item(3,3).Value = 30.1        # C3
$wsInterannual.Cells.Item(3,4).Value = 23          # D3
$wsInterannual.Cells.Item(3,5).Value = 24.2        # E3

# ---------------------------------------------------------------------------
# 5. Add a brand-new sheet "Major update - High Priority " at the end that
#    keeps the original ("major update") High Priority break-up figures.
# ---------------------------------------------------------------------------
$wsMajor = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMajor.Name = "Major update - High Priority "

$wsMajor.Cells.Item(1,1).Value = "Break-up"
$wsMajor.Cells.Item(1,2).Value = "High Species (no.)"
$wsMajor.Cells.Item(1,3).Value = "High Species (perc.)"
$wsMajor.Cells.Item(1,4).Value = "New High Species (no.)"
$wsMajor.Cells.Item(1,5).Value = "New High Species (perc.)"
$wsMajor.Range("A1:E1").Font.Bold = $true
$wsMajor.Range("A1:E1").HorizontalAlignment = -4108

$wsMajor.Cells.Item(2,1).Value = "Trend New"
$wsMajor.Cells.Item(2,2).Value = 3
$wsMajor.Cells.Item(2,3).Value = 23.1
$wsMajor.Cells.Item(2,4).Value = 3
$wsMajor.Cells.Item(2,5).Value = 23.1

$wsMajor.Cells.Item(3,1).Value = "IUCN"
$wsMajor.Cells.Item(3,2).Value = 10
$wsMajor.Cells.Item(3,3).Value = 76.90000000000001
$wsMajor.Cells.Item(3,4).Value = 10
$wsMajor.Cells.Item(3,5).Value = 76.90000000000001
